# Generate Report for Handoff
# Adds two new tracked files (a .md dependency source and two .png files
# that depend on it) to the localization-status workbook: one summary row
# per file on "Overview", and one detail row per file on each locale sheet
# ("zh-cn" / "de-de"). Also refreshes the existing row's "Latest Handoff"
# timestamps to reflect the new handoff pass.
#
# Hyperlink objects in this engine can only be appended (re-assigning
# .Address/.TextToDisplay on an existing Hyperlink, or re-Add()-ing one on
# top of a cell that already has one, leaves the stale entry behind), so
# every sheet's Hyperlinks collection is cleared and rebuilt from scratch,
# in column-major / row-major reading order, which reproduces the rId2,
# rId3, rId4 ... numbering a fresh build would get.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$status  = "Ready for handoff"
$commit  = "4dec7402d9945e50af442a1380d2764dc9d6606b"

# New source / dependent files introduced by this handoff.
$pngFile1 = "2a3b59e6-55b6-4ddc-be6a-c2fb1f22f009.png"
$mdFile   = "3a0cddb6-bc0b-4c20-8931-d9444b9d40d3.md"
$pngFile2 = "76b38c36-8554-4ec4-83db-86c3ba30958d.png"

$overviewDate = "2016-41-12 06:41:48"

$srcUrl1 = "https://github.com/OpenLocalizationTest/oltest/blob/$commit/e2e/$pngFile1"
$srcUrlMd = "https://github.com/OpenLocalizationTest/oltest/blob/$commit/e2e/$mdFile"
$srcUrl2 = "https://github.com/OpenLocalizationTest/oltest/blob/$commit/e2e/$pngFile2"

# ---------------------------------------------------------------------
# Overview sheet: update the existing row's handoff date, then append one
# summary row for each of the two new files.
# ---------------------------------------------------------------------
$overview.Range("D2").Value = $overviewDate

$overview.Range("B3").Value = $status
$overview.Range("C3").Value = $status
$overview.Range("D3").Value = $overviewDate

$overview.Range("B4").Value = $status
$overview.Range("C4").Value = $status
$overview.Range("D4").Value = $overviewDate

$overview.Hyperlinks.Delete()
$overview.Hyperlinks.Add($overview.Range("A2"), $srcUrl1, "", "", $pngFile1)
$overview.Hyperlinks.Add($overview.Range("A3"), $srcUrlMd, "", "", $mdFile)
$overview.Hyperlinks.Add($overview.Range("A4"), $srcUrl2, "", "", $pngFile2)

Write-Output "Overview sheet updated"

# ---------------------------------------------------------------------
# Per-locale detail sheets (zh-cn, de-de): row 2 is repointed from the old
# .md/.xlf pair to the new first .png; row 3 carries the new .md row (the
# dependency source, handed off with reason "Include"); row 4 carries the
# second new .png row (handed off with reason "IsDependency" on the .md).
# ---------------------------------------------------------------------
function Update-LocaleSheet {
    param(
        $ws,
        [string]$locale,
        [string]$xlfHash,
        [string]$handoffDatetime,
        [string]$olhandoffCommit
    )

    $zhHandback = "3a0cddb6-bc0b-4c20-8931-d9444b9d40d3.$xlfHash.$locale.xlf"
    $dTarget1 = "1a3eab3540262c430ecda44078916c54611daca4.png"
    $dTarget2 = "50b6a3b267138c6e74bfbaa5d3bcac0e028e1a2d.png"

    $handbackUrlBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$olhandoffCommit/ol-handoff/OpenLocalizationTestOrg/oltest.$locale/ci/ht"

    # Row 2: 1st new .png -- a dependency of the .md (handed off because of
    # its dependency, no direct target file of its own yet).
    $ws.Range("C2").Value = $status
    $ws.Range("H2").Value = "0001-01-01 00:00:00"
    $ws.Range("I2").Value = "IsDependency"
    $ws.Range("J2").Value = "e2e\$mdFile"
    $ws.Range("E2").Value = $handoffDatetime
    $ws.Range("E2").NumberFormat = "yyyy-mm-dd hh:mm:ss"

    # Row 3: the .md itself -- Included, and produced the real .xlf handback.
    $ws.Range("C3").Value = $status
    $ws.Range("H3").Value = "0001-01-01 00:00:00"
    $ws.Range("I3").Value = "Include"
    $ws.Range("E3").Value = $handoffDatetime
    $ws.Range("E3").NumberFormat = "yyyy-mm-dd hh:mm:ss"

    # Row 4: 2nd new .png -- same shape as row 2's dependency row.
    $ws.Range("C4").Value = $status
    $ws.Range("H4").Value = "0001-01-01 00:00:00"
    $ws.Range("I4").Value = "IsDependency"
    $ws.Range("J4").Value = "e2e\$mdFile"
    $ws.Range("E4").Value = $handoffDatetime
    $ws.Range("E4").NumberFormat = "yyyy-mm-dd hh:mm:ss"

    $ws.Hyperlinks.Delete()

    $ws.Hyperlinks.Add($ws.Range("A2"), $srcUrl1, "", "", $pngFile1)
    $ws.Hyperlinks.Add($ws.Range("B2"), $srcUrl1, "", "", ".png")
    $ws.Hyperlinks.Add($ws.Range("D2"), "$handbackUrlBase/$dTarget1", "", "", $dTarget1)

    $ws.Hyperlinks.Add($ws.Range("A3"), $srcUrlMd, "", "", $mdFile)
    $ws.Hyperlinks.Add($ws.Range("B3"), $srcUrlMd, "", "", ".md")
    $ws.Hyperlinks.Add($ws.Range("D3"), "$handbackUrlBase/$zhHandback", "", "", $zhHandback)

    $ws.Hyperlinks.Add($ws.Range("A4"), $srcUrl2, "", "", $pngFile2)
    $ws.Hyperlinks.Add($ws.Range("B4"), $srcUrl2, "", "", ".png")
    $ws.Hyperlinks.Add($ws.Range("D4"), "$handbackUrlBase/$dTarget2", "", "", $dTarget2)
}

Update-LocaleSheet $zhcn "zh-cn" "c72712d0d988816f0e7d8e27288f55cf361ab8f6" "2016-03-12 06:41:45" "bd69457f97bcdf6da9602e0fefa762b0884aca71"
Write-Output "zh-cn sheet updated"

Update-LocaleSheet $dede "de-de" "c72712d0d988816f0e7d8e27288f55cf361ab8f6" "2016-03-12 06:41:48" "49a515bec0b1a405abb2bee1d640ef379bda7d8a"
Write-Output "de-de sheet updated"

Write-Output "Done"
